$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This document's paragraphs keep their formatting (pPr/rPr) in place while
# several blocks of text content get relocated to different paragraphs/runs.
# Strategy: capture every original text value we need first (before any
# mutation), then write the new values into their destination locations.
# ---------------------------------------------------------------------------

$docEnd = $d.Content.End

# --- Step 1: locate & capture the three labelled values inside the
#     "Avaliação" bullet paragraph (Paragraphs(17)) using Find, while
#     offsets are still untouched. ------------------------------------------
$rMetodoLabel = $d.Range(0, $docEnd)
$null = $rMetodoLabel.Find.Execute("Método: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rCriterioLabel = $d.Range($rMetodoLabel.End, $docEnd)
$null = $rCriterioLabel.Find.Execute("Critério: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rNormaLabel = $d.Range($rCriterioLabel.End, $docEnd)
$null = $rNormaLabel.Find.Execute("Norma de recuperação: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Paragraphs(17).Range.End points one position past the last visible
# character (it includes the trailing paragraph-mark/CR), so back it off by
# one to get a range that ends right after the last visible character.
$p17End = $d.Paragraphs(17).Range.End - 1

$valMetodoRange = $d.Range($rMetodoLabel.End, $rCriterioLabel.Start)
$valCriterioRange = $d.Range($rCriterioLabel.End, $rNormaLabel.Start)
$valNormaRange = $d.Range($rNormaLabel.End, $p17End)

$nl = [char]11

# The "Método:" and "Critério:" value runs each end with a <w:br/> (chr 11)
# that visually separates them from the next bold label; that break is part
# of the *paragraph-17 slot structure*, not of the text payload itself, so it
# must stay behind at the slot rather than travel with the relocated text.
$origValMetodo = $valMetodoRange.Text
if ($origValMetodo.Substring($origValMetodo.Length - 1, 1) -eq $nl) {
    $origValMetodo = $origValMetodo.Substring(0, $origValMetodo.Length - 1)
}
$origValCriterio = $valCriterioRange.Text
if ($origValCriterio.Substring($origValCriterio.Length - 1, 1) -eq $nl) {
    $origValCriterio = $origValCriterio.Substring(0, $origValCriterio.Length - 1)
}
$origValNorma = $valNormaRange.Text
if ($origValNorma.Substring($origValNorma.Length - 1, 1) -eq $nl) {
    $origValNorma = $origValNorma.Substring(0, $origValNorma.Length - 1)
}

# --- Step 2: capture the original whole-paragraph texts that will move. ---
# NOTE: Paragraphs(N).Range.Text includes a trailing paragraph-mark (CR,
# chr 13) character. Carrying that CR into a different paragraph's Range.Text
# would insert a spurious paragraph break, so it is stripped immediately
# after capture.
$origP6 = $d.Paragraphs(6).Range.Text    # "Fornecer aos alunos..." (PT objetivo)
$origP6 = $origP6.Substring(0, $origP6.Length - 1)
$origP7 = $d.Paragraphs(7).Range.Text    # "Supply the students..." (EN objetivo, italic)
$origP7 = $origP7.Substring(0, $origP7.Length - 1)
$origP9 = $d.Paragraphs(9).Range.Text    # "7455355 - Robson da Silva Rocha"
$origP9 = $origP9.Substring(0, $origP9.Length - 1)
$origP11 = $d.Paragraphs(11).Range.Text  # "Tecnologias de Tratamento..." (PT resumo)
$origP11 = $origP11.Substring(0, $origP11.Length - 1)
$origP12 = $d.Paragraphs(12).Range.Text  # "Water Treatment Technologies..." (EN resumo, italic)
$origP12 = $origP12.Substring(0, $origP12.Length - 1)
$origP14 = $d.Paragraphs(14).Range.Text  # "- Características das águas..." (PT programa)
$origP14 = $origP14.Substring(0, $origP14.Length - 1)
$origP19 = $d.Paragraphs(19).Range.Text  # Bibliografia reference text
$origP19 = $origP19.Substring(0, $origP19.Length - 1)

# --- Step 3: write the relocated text into the "Avaliação" sub-ranges. ----
# (Do this before touching the other paragraphs so the captured Find-based
# ranges are still valid / unaffected by upstream length changes.)
$valNormaRange.Text = $origP19
$valCriterioRange.Text = $origValNorma
$valMetodoRange.Text = $origValCriterio

# --- Step 4: write the relocated text into the other paragraphs. ----------
# Cycle: P6 -> P9 -> P19 -> (Norma value, handled above) ; P11 -> P6 ;
#        P14 -> P11 ; (Método value, handled above) -> P14
# Swap:  P7 <-> P12
$d.Paragraphs(19).Range.Text = $origP9
$d.Paragraphs(9).Range.Text = $origP6
$d.Paragraphs(6).Range.Text = $origP11
$d.Paragraphs(11).Range.Text = $origP14
$d.Paragraphs(14).Range.Text = $origValMetodo

$d.Paragraphs(7).Range.Text = $origP12
$d.Paragraphs(12).Range.Text = $origP7
